{"js": "// V1.1 avec gestion Autre type Intervention\n//\n// Each \"{messageInfoXxx}\" placeholder paragraph used to be prefixed with a\n// red \"NB : \" (or \"NB : \") label. The label is removed, leaving the\n// paragraph starting directly at \"{\".\n//\n// We search with a wildcard pattern that matches \"NB\" followed by any\n// characters up to (and including) the opening \"{\" of the placeholder,\n// then replace the whole matched range with just \"{\" - this removes the\n// \"NB\" + separator runs while leaving the \"{...}\" placeholder runs (and\n// their formatting/spell-check markers) completely untouched.\nconst results = context.document.body.search(\"NB*{\", { matchWildcards: true });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"{\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# V1.1 avec gestion Autre type Intervention\n# Remove the leading \"NB : \" label (plain or non-breaking space variants)\n# that precedes each \"{messageInfoXxx}\" placeholder paragraph.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Text = \"NB : \"\n$rng.Find.MatchCase = $false\n$rng.Find.MatchWholeWord = $false\n$rng.Find.MatchWildcards = $false\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 1\n$rng.Find.Replacement.Text = \"\"\n$rng.Find.Execute(\n  [ref]\"NB : \",\n  [ref]$false,\n  [ref]$false,\n  [ref]$false,\n  [ref]$false,\n  [ref]$false,\n  [ref]$true,\n  [ref]1,\n  [ref]$false,\n  [ref]\"\",\n  [ref]2\n) | Out-Null\n"}
